$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.969.60"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "1.743.55"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "311.53"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "0.4971"
$ws.Range("E7").Value = "  +7.17%  "
$ws.Range("D8").Value = "0.3578"
$ws.Range("E8").Value = "  +3.61%  "
$ws.Range("D9").Value = "42.31"
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("D10").Value = "0.07269"
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("D11").Value = "1.061"
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "20.22"
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").Value = "5.957"
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("D15").Value = "1.743.25"
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("D16").Value = "6.851"
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("D17").Value = "87.16"
$ws.Range("E17").Value = "  -2.97%  "
$ws.Range("D18").Value = "0.00001038"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("D19").Value = "0.06367"
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "16.59"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Value = "5.726"
$ws.Range("E22").Value = "  +1.25%  "
$ws.Range("D23").Value = "27.037.59"
$ws.Range("E23").Value = "  -0.62%  "
$ws.Range("D24").Value = "11.33"
$ws.Range("E24").Value = "  +4.36%  "
$ws.Range("E25").Value = "  -4.20%  "
$ws.Range("D26").Value = "156.10"
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("D27").Value = "19.85"
$ws.Range("E27").Value = "  +1.48%  "
$ws.Range("D28").Value = "1.943.83"
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("D29").Value = "2.195"
$ws.Range("E29").Value = "  +1.54%  "
$ws.Range("D30").Value = "120.57"
$ws.Range("E30").Value = "  +0.93%  "
$ws.Range("D31").Value = "1.051"
$ws.Range("E31").Value = "  +1.34%  "
$ws.Range("D32").Value = "0.09509"
$ws.Range("E32").Value = "  +4.19%  "
$ws.Range("D33").Value = "3.579"
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("D34").Value = "5.378"
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("D35").Value = "0.02203"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").Value = "0.05875"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "11.05"
$ws.Range("E37").Value = "  -1.22%  "
$ws.Range("E38").Value = "  +2.54%  "
$ws.Range("D39").Value = "0.2001"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Value = "4.779"
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").Value = "0.6025"
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("D42").Value = "1.111"
$ws.Range("E42").Value = "  -1.96%  "
$ws.Range("D43").Value = "7.597"
$ws.Range("E43").Value = "  +0.67%  "
$ws.Range("D44").Value = "12.82"
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("D45").Value = "3.597"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("D46").Value = "0.5668"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("D47").Value = "120.44"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("D48").Value = "1.860"
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("D49").Value = "0.06678"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").Value = "1.100"
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("E51").Value = "  +0.15%  "
